$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "X-n251-k28"
$ws.Range("A17").Value = "X-n200-k36"
$ws.Range("D17").Value = "57666,2…"
$ws.Range("A18").Value = "X-n228-k23"

$ws.Range("B17").Value = 58578
$ws.Range("C17").Value = 83

$ws.Range("B18").Value = 25742
$ws.Range("C18").Value = 156

$ws.Range("B19").Value = 38684
$ws.Range("C19").Value = 16

$ws.Range("A20").Value = "X-n256-k16"
$ws.Range("C20").Value = 27

$ws.Range("C19").Select()
